# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2210
#   *_new  -> *_FV2304
# Then wrap the data range in an Excel Table and freeze the header row,
# matching the authored commit "chore: adapt column header formatting to
# respective input file names (#7)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the 21 column headers in row 1 (A1:U1) to use the new
#    "<formatversion>" suffixes instead of the old "_old"/"_new" ones.
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# 2. Turn the used range A1:U88 into a real Excel Table ("Table1") so the
#    headers double as filter buttons / structured references.
$dataRange = $ws.Range("A1:U88")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# 3. Freeze the header row (split below row 1, top-left of the scrolling
#    area is A2) so the header stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
